# Update the cryptocurrency price table with the latest scraped values.
# A leading apostrophe forces Excel to store a numeric-looking string
# (e.g. "227.70", "0.999") as literal text instead of a Number,
# which preserves trailing zeros / exact formatting shown on the site.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '34.620.47'
$ws.Range("E2").Value = '  +1.28%  '

$ws.Range("D3").Value = '1.802.03'
$ws.Range("E3").Value = '  +1.22%  '

$ws.Range("D4").Value = '''0.999'
$ws.Range("E4").Value = '  -0.34%  '

$ws.Range("D5").Value = '''227.70'
$ws.Range("E5").Value = '  +0.73%  '

$ws.Range("E6").Value = '  +1.88%  '

$ws.Range("E7").Value = '  -0.37%  '

$ws.Range("D8").Value = '''32.80'
$ws.Range("E8").Value = '  +3.40%  '

$ws.Range("E9").Value = '  +1.90%  '

$ws.Range("E10").Value = '  +1.01%  '

$ws.Range("E11").Value = '  +0.34%  '

$ws.Range("D12").Value = '2.061.76'
$ws.Range("E12").Value = '  +1.02%  '

$ws.Range("D13").Value = '''11.15'
$ws.Range("E13").Value = '  +1.97%  '

$ws.Range("D14").Value = '1.794.11'
$ws.Range("E14").Value = '  +1.30%  '

$ws.Range("D15").Value = '''0.640'
$ws.Range("E15").Value = '  +2.77%  '

$ws.Range("D16").Value = '34.599.95'
$ws.Range("E16").Value = '  +1.28%  '

$ws.Range("E17").Value = '  +3.93%  '

$ws.Range("D18").Value = '''68.90'
$ws.Range("E18").Value = '  +1.58%  '

$ws.Range("D19").Value = '0.0₃0805'

$ws.Range("D20").Value = '''247.46'
$ws.Range("E20").Value = '  +0.65%  '

$ws.Range("D21").Value = '''11.36'
$ws.Range("E21").Value = '  +3.75%  '

$ws.Range("D22").Value = '''0.999'
$ws.Range("E22").Value = '  -0.30%  '

$ws.Range("E23").Value = '  +2.84%  '

$ws.Range("D24").Value = '''169.80'
$ws.Range("E24").Value = '  +4.74%  '

$ws.Range("E25").Value = '  +2.05%  '

$ws.Range("D26").Value = '''7.32'
$ws.Range("E26").Value = '  +1.94%  '

$ws.Range("E27").Value = '  +2.06%  '

$ws.Range("D28").Value = '''0.117'
$ws.Range("E28").Value = '  +2.53%  '

$ws.Range("E29").Value = '  -0.46%  '

$ws.Range("D30").Value = '''4.13'
$ws.Range("E30").Value = '  +11.10%  '

$ws.Range("B31").Value = 'PancakeSwap'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D31").Value = '''1.25'
$ws.Range("E31").Value = '  +0.98%  '

$ws.Range("B32").Value = 'Hedera'
$ws.Range("C32").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D32").Value = '''0.0527'
$ws.Range("E32").Value = '  +1.39%  '

$ws.Range("B33").Value = 'Filecoin'
$ws.Range("C33").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D33").Value = '''3.82'
$ws.Range("E33").Value = '  +2.40%  '

$ws.Range("E34").Value = '  +3.10%  '

$ws.Range("D35").Value = '1.432.24'
$ws.Range("E35").Value = '  -0.51%  '

$ws.Range("D36").Value = '''2.60'
$ws.Range("E36").Value = '  +8.66%  '

$ws.Range("D37").Value = '''0.677'
$ws.Range("E37").Value = '  +3.74%  '

$ws.Range("E38").Value = '  +3.12%  '

$ws.Range("E39").Value = '  +0.55%  '

$ws.Range("D40").Value = '''85.29'
$ws.Range("E40").Value = '  +6.49%  '

$ws.Range("E41").Value = '  +2.96%  '

$ws.Range("E42").Value = '  +2.20%  '

$ws.Range("E43").Value = '  +3.55%  '

$ws.Range("D44").Value = '''13.84'
$ws.Range("E44").Value = '  +2.57%  '

$ws.Range("E45").Value = '  +3.06%  '

$ws.Range("E46").Value = '  +0.75%  '

$ws.Range("E47").Value = '  +0.61%  '

$ws.Range("D48").Value = '1.960.78'
$ws.Range("E48").Value = '  +0.90%  '

$ws.Range("D49").Value = '''105.83'
$ws.Range("E49").Value = '  +1.45%  '

$ws.Range("D50").Value = '''0.999'
$ws.Range("E50").Value = '  -0.32%  '

$ws.Range("E51").Value = '  -4.51%  '
